# Produtos.xlsx — "Pegando o preço de um produto"
# Swap the tracked product from an iPhone 11 / listing-search links to an
# iPhone 15 (with direct product-page links), and refresh the Smart TV's
# retailer links to their current product pages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New product links (direct product pages instead of old search/listing URLs).
$amazonIphone = "https://www.amazon.com.br/gp/aw/d/B0CP6CVJSG/?_encoding=UTF8&pd_rd_plhdr=t&aaxitk=71870feb0a4de624959d0dc6cd3143de&hsa_cr_id=2734213800601&qid=1751583756&sr=1-2-9e67e56a-6f64-441f-a281-df67fc737124&ref_=sbx_be_s_sparkle_dlcd_asin_1_img&pd_rd_w=2UCbV&content-id=amzn1.sym.da990715-51de-4373-b8ab-d245f18f3edf%3Aamzn1.sym.da990715-51de-4373-b8ab-d245f18f3edf&pf_rd_p=da990715-51de-4373-b8ab-d245f18f3edf&pf_rd_r=D3D8M9T0W18MAZZJM6Z8&pd_rd_wg=EvWpu&pd_rd_r=49a8f4fa-0d1f-47d6-afe7-8a6c94744c0c&th=1"
$americanasIphone = "https://www.americanas.com.br/apple-iphone-15-de-128gb-preto-174n0p87924z7261/p?idsku=7763713&utm_source=YSMESP&utm_medium=buscappc&utm_campaign=celulares_3p&utm_content=bp_pl_sh_go_aloc_digital_celulares_topitens_3p_aon25-00014&utm_term=pla_shoppingpadrao&gad_source=1&gad_campaignid=22388929149&gbraid=0AAAAAD37VppWcVmTKgnAlwmq47CtWjCxO&gclid=Cj0KCQjw1JjDBhDjARIsABlM2Ssnjp0NqiipnAJnKBFoSQZiWSwxJ4lUhNNkWkhXS1dUGeHBf25pXf8aAhrvEALw_wcB"
$magazineIphone = "https://www.magazineluiza.com.br/apple-iphone-15-128gb-preto-61-48mp-ios-5g/p/238035600/te/ip15/?&seller_id=magazineluiza&utm_source=google&utm_medium=cpc&utm_term=79684&utm_campaign=google_eco_per_ven_pla_mob_apo_1p_te-csp&utm_content=&partner_id=79684&gclsrc=aw.ds&gad_source=1&gad_campaignid=22578732921&gbraid=0AAAAAD4zZmTvhM3LLUO76XwjDkjAbe4ca&gclid=Cj0KCQjw1JjDBhDjARIsABlM2SsbLJ7_-kwUaY-DHsgDCKDYJ50rBTUn43KmU8TLEQ8Hvp-T-P-oSAgaAn98EALw_wcB"

$amazonTv = "https://www.amazon.com.br/LG-50NANO80T-Processador-Chromecast-integrado/dp/B0D3JB8GYZ/ref=asc_df_B0D3JB8GYZ?mcid=769db38cdb513c40a442f4c7cc227410&tag=googleshopp00-20&linkCode=df0&hvadid=709964502896&hvpos=&hvnetw=g&hvrand=17184299404204777997&hvpone=&hvptwo=&hvqmt=&hvdev=c&hvdvcmdl=&hvlocint=&hvlocphy=1031811&hvtargid=pla-2323939420537&psc=1&language=pt_BR&gad_source=1"
$americanasTv = "https://www.americanas.com.br/smart-tv-4k-50-polegadas-lg-uhd-50ut8050-processador-a5-ger7-ai-alexa-chromecast-integrado-otimizador-de-jogos-webos-24-controle-smart-magic-7507042514/p?idsku=5725301&utm_source=YSMESP&utm_medium=buscappc&utm_campaign=alwayson-25&utm_content=bp_pl_sh_go_aloc_digital_apostassortimento_3p_aon25-00234&utm_term=pla_shopping&gad_source=1&gad_campaignid=22744334395&gbraid=0AAAAAD37Vpp_JV01DVyUO7oovGs82Yucz&gclid=Cj0KCQjw1JjDBhDjARIsABlM2SunsoJy03aCCv1IHpsRxNNoRQqvqLxcybiTTcbs3YwbiicYJS3uZWsaAmpdEALw_wcB"
$magazineTv = "https://www.magazineluiza.com.br/smart-tv-50-4k-lg-ultra-hd-50ut8050-com-processador-a5-ger7-ai-alexa-chromecast-integrado-otimizador-de-jogos-webos-24-e-controle-smart-magic/p/238617500/et/tves/?&seller_id=magazineluiza&utm_source=google&utm_medium=cpc&utm_term=79710&utm_campaign=google_eco_per_ven_pla_tc_apo_1p_et-csp&utm_content=&partner_id=79710&gclsrc=aw.ds&gad_source=1&gad_campaignid=22561248932&gbraid=0AAAAAD4zZmTItqqzq1KWwBaXEdlBOIQCI&gclid=Cj0KCQjw1JjDBhDjARIsABlM2SvijLU9hiNk_RLA5c6s1AMxr70tz3SwbtfO8u5P_6qLd55AzURost4aAiJtEALw_wcB"

# Drop every existing hyperlink so we can rebuild them cleanly (Add() keeps
# stacking new relationships on top of old ones if they aren't cleared).
$ws.Hyperlinks.Delete()

# Row 2: iPhone 11 -> iPhone 15
$ws.Range("A2").Value2 = "iPhone 15 Apple 128GB Preto"

$ws.Hyperlinks.Add($ws.Range("B2"), $amazonIphone, [Type]::Missing, [Type]::Missing, $amazonIphone) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), $americanasIphone, [Type]::Missing, [Type]::Missing, $americanasIphone) | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), $magazineIphone, [Type]::Missing, [Type]::Missing, $magazineIphone) | Out-Null

# Row 3: Smart TV links refreshed to current product pages
$ws.Hyperlinks.Add($ws.Range("B3"), $amazonTv, [Type]::Missing, [Type]::Missing, $amazonTv) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), $americanasTv, [Type]::Missing, [Type]::Missing, $americanasTv) | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), $magazineTv, [Type]::Missing, [Type]::Missing, $magazineTv) | Out-Null

# Hyperlinks.Add() re-styles the touched cells with a freshly minted style;
# put them back on the workbook's original "Hiperlink" cell style.
$ws.Range("B2:D3").Style = "Hiperlink"

# Rows 2 & 3 grew slightly taller in the resaved file.
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15

# Resaved file's selection landed on B3.
$ws.Range("B3").Select() | Out-Null
